# Actualización automatica mar abr  6 17:33:39 CEST 2021
#
# The dataset gained a new "Extranjeros" (misspelled "Extrenjeros" in the
# source data, including the lower-case slug and the iaest-measure id) column.
# It is inserted as the new column C, pushing every existing column from the
# former C onward one position to the right (C->D, D->E, ... L->M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this automatically shifts columns C..L (and all
# their row data / formatting) one place to the right, becoming D..M.
$ws.Columns("C").Insert()

# Populate the newly freed column C with the "Extrenjeros" metadata, mirroring
# the same five metadata rows used by every other column in the sheet.
$ws.Range("C1").Value = "Extrenjeros"
$ws.Range("C2").Value = "extrenjeros"
$ws.Range("C3").Value = "iaest-measure:extrenjeros"
$ws.Range("C4").Value = "medida"
$ws.Range("C5").Value = "xsd:int"
